$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers; force Text format
# so Excel keeps them as literal strings (matching the source data which
# uses these as formatted text, not numeric values).
$textForceAddrs = @("D5", "D6", "D7", "D10", "D11", "D16", "D19", "D21", "D22", "D23", "D24", "D32", "D33", "D35", "D36", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range('D5').Value = '484.59'
$ws.Range('D6').Value = '143.66'
$ws.Range('D7').Value = '0.999'
$ws.Range('D10').Value = '5.76'
$ws.Range('D11').Value = '0.0961'
$ws.Range('D16').Value = '20.85'
$ws.Range('D19').Value = '4.49'
$ws.Range('D21').Value = '314.93'
$ws.Range('D22').Value = '0.997'
$ws.Range('D23').Value = '5.79'
$ws.Range('D24').Value = '58.23'
$ws.Range('D32').Value = '147.60'
$ws.Range('D33').Value = '18.10'
$ws.Range('D35').Value = '5.12'
$ws.Range('D36').Value = '1.14'
$ws.Range('D39').Value = '33.70'
$ws.Range('D41').Value = '0.996'
$ws.Range('D42').Value = '0.0547'
$ws.Range('D43').Value = '0.596'
$ws.Range('D45').Value = '259.20'
$ws.Range('D46').Value = '0.0920'
$ws.Range('D47').Value = '10.16'
$ws.Range('D48').Value = '4.66'
$ws.Range('D50').Value = '17.35'
$ws.Range('D2').Value = '55.822.48'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '2.444.73'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +3.67%  '
$ws.Range('E6').Value = '  +9.95%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +2.53%  '
$ws.Range('D9').Value = '2.445.00'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  +8.75%  '
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('E12').Value = '  +4.75%  '
$ws.Range('E13').Value = '  +1.39%  '
$ws.Range('D14').Value = '2.876.20'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').Value = '55.864.87'
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('E16').Value = '  +6.44%  '
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '2.451.40'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('E19').Value = '  +6.32%  '
$ws.Range('E20').Value = '  +4.78%  '
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  +7.33%  '
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('E25').Value = '  +5.47%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('E27').Value = '  +2.69%  '
$ws.Range('D28').Value = '2.564.88'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  +7.15%  '
$ws.Range('D30').Value = '0.0₃0770'
$ws.Range('E30').Value = '  +8.41%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  +8.21%  '
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('E38').Value = '  +6.39%  '
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('E40').Value = '  +7.52%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('E42').Value = '  +4.56%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +6.50%  '
$ws.Range('E45').Value = '  +11.93%  '
$ws.Range('E46').Value = '  +4.15%  '
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('E48').Value = '  +11.80%  '
$ws.Range('E49').Value = '  +4.55%  '
$ws.Range('E50').Value = '  +4.76%  '
$ws.Range('D51').Value = '1.855.58'
$ws.Range('E51').Value = '  -3.82%  '
